$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1140
$ws.Range("I43").Value = 1140
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 1140
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -1071
$ws.Range("N43").ClearContents()

$ws.Range("H129").Value = 924.8182
$ws.Range("J129").Value = 1009.7083
$ws.Range("L129").Value = 3029.1249
$ws.Range("N129").Value = -13029.1249

$ws.Range("H138").Value = 4931.4614
$ws.Range("I138").Value = 1693.0476
$ws.Range("J138").Value = 7125.2256
$ws.Range("K138").Value = 5079.142800000001
$ws.Range("L138").Value = 21375.6768
$ws.Range("M138").Value = 60.85719999999947
$ws.Range("N138").Value = -31655.6768

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 39799.49
$ws.Range("I32").Value = 7739.1724
$ws.Range("K32").Value = 7739.1724
$ws.Range("M32").Value = -7452.1724

$ws.Range("H45").Value = 1541.6875
$ws.Range("I45").Value = 1681.6666
$ws.Range("J45").Value = 1457.7
$ws.Range("K45").Value = 1681.6666
$ws.Range("L45").Value = 1457.7
$ws.Range("M45").Value = -1304.6666
$ws.Range("N45").Value = -2211.7

$ws.Range("H63").Value = 2428.5
$ws.Range("I63").Value = 1021.25
$ws.Range("J63").Value = 3366.6667
$ws.Range("K63").Value = 1021.25
$ws.Range("L63").Value = 3366.6667
$ws.Range("M63").Value = -335.25
$ws.Range("N63").Value = -4738.6667

$ws.Range("H66").Value = 2428.5
$ws.Range("I66").Value = 1021.25
$ws.Range("J66").Value = 3366.6667
$ws.Range("K66").Value = 5106.25
$ws.Range("L66").Value = 16833.3335
$ws.Range("M66").Value = -1674.25
$ws.Range("N66").Value = -23697.3335

$ws.Range("H88").Value = 8301.200000000001
$ws.Range("I88").Value = 1335.3334
$ws.Range("J88").Value = 18750
$ws.Range("K88").Value = 1335.3334
$ws.Range("L88").Value = 18750
$ws.Range("M88").Value = -929.3334
$ws.Range("N88").Value = -19562

$ws.Range("H91").Value = 8301.200000000001
$ws.Range("I91").Value = 1335.3334
$ws.Range("J91").Value = 18750
$ws.Range("K91").Value = 1335.3334
$ws.Range("L91").Value = 18750
$ws.Range("M91").Value = 68.66660000000002
$ws.Range("N91").Value = -21558

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 86937.08
$ws.Range("I86").Value = 139024.12
$ws.Range("J86").Value = 3597.8
$ws.Range("K86").Value = 139024.12
$ws.Range("L86").Value = 3597.8
$ws.Range("M86").Value = -137901.12
$ws.Range("N86").Value = -5843.8

$ws.Range("H89").Value = 86937.08
$ws.Range("I89").Value = 139024.12
$ws.Range("J89").Value = 3597.8
$ws.Range("K89").Value = 695120.6
$ws.Range("L89").Value = 17989
$ws.Range("M89").Value = -689504.6
$ws.Range("N89").Value = -29221

$ws.Range("H105").Value = 224298.89
$ws.Range("I105").Value = 168279.83
$ws.Range("K105").Value = 168279.83
$ws.Range("M105").Value = -166532.83

$ws.Range("H107").Value = 18520636
$ws.Range("I107").Value = 30304042
$ws.Range("J107").Value = 3858.8572
$ws.Range("K107").Value = 30304042
$ws.Range("L107").Value = 3858.8572
$ws.Range("M107").Value = -30302122
$ws.Range("N107").Value = -7698.8572

$ws.Range("H130").Value = 45157.25
$ws.Range("I130").Value = 30709
$ws.Range("J130").Value = 49973.332
$ws.Range("K130").Value = 30709
$ws.Range("L130").Value = 49973.332
$ws.Range("M130").Value = -25689
$ws.Range("N130").Value = -60013.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15516.634
$ws.Range("I31").Value = 25934.324
$ws.Range("J31").Value = 2074.4517
$ws.Range("K31").Value = 25934.324
$ws.Range("L31").Value = 2074.4517
$ws.Range("M31").Value = -25639.324
$ws.Range("N31").Value = -2664.4517

$ws.Range("H34").Value = 15516.634
$ws.Range("I34").Value = 25934.324
$ws.Range("J34").Value = 2074.4517
$ws.Range("K34").Value = 25934.324
$ws.Range("L34").Value = 2074.4517
$ws.Range("M34").Value = -25732.324
$ws.Range("N34").Value = -2478.4517

$ws.Range("H50").Value = 11095
$ws.Range("J50").Value = 11095
$ws.Range("L50").Value = 11095
$ws.Range("N50").Value = -12345

$ws.Range("H51").Value = 6862.5713
$ws.Range("J51").Value = 7979.6
$ws.Range("L51").Value = 7979.6
$ws.Range("N51").Value = -9451.6

$ws.Range("H60").Value = 18073.334
$ws.Range("J60").Value = 19288
$ws.Range("L60").Value = 19288
$ws.Range("N60").Value = -20310

$ws.Range("H61").Value = 6862.5713
$ws.Range("J61").Value = 7979.6
$ws.Range("L61").Value = 7979.6
$ws.Range("N61").Value = -8675.6

$ws.Range("H62").Value = 3657
$ws.Range("I62").Value = 2360.8
$ws.Range("K62").Value = 2360.8
$ws.Range("M62").Value = -1736.8

$ws.Range("H65").Value = 3657
$ws.Range("I65").Value = 2360.8
$ws.Range("K65").Value = 12505
$ws.Range("M65").Value = -8684

$ws.Range("H74").Value = 19181.084
$ws.Range("J74").Value = 19181.084
$ws.Range("L74").Value = 19181.084
$ws.Range("N74").Value = -20929.084

$ws.Range("H77").Value = 19181.084
$ws.Range("J77").Value = 19181.084
$ws.Range("L77").Value = 57543.25199999999
$ws.Range("N77").Value = -66279.25199999999

$ws.Range("H122").Value = 2214.8147
$ws.Range("I122").Value = 2532
$ws.Range("J122").Value = 1818.3334
$ws.Range("K122").Value = 7596
$ws.Range("L122").Value = 5455.0002
$ws.Range("M122").Value = -5146
$ws.Range("N122").Value = -10355.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 584.0476
$ws.Range("I114").Value = 130.77777
$ws.Range("J114").Value = 924
$ws.Range("K114").Value = 392.33331
$ws.Range("L114").Value = 2772
$ws.Range("M114").Value = 2861.66669
$ws.Range("N114").Value = -9280

$ws.Range("H131").Value = 4025.5557
$ws.Range("J131").Value = 5000
$ws.Range("L131").Value = 15000
$ws.Range("N131").Value = -25080

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()

$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

$ws.Range("H70").Value = 66765.97
$ws.Range("I70").Value = 94966.55
$ws.Range("K70").Value = 94966.55
$ws.Range("M70").Value = -94696.55

$ws.Range("H73").Value = 66765.97
$ws.Range("I73").Value = 94966.55
$ws.Range("K73").Value = 94966.55
$ws.Range("M73").Value = -94030.55

$ws.Range("H132").Value = 2930.8696
$ws.Range("I132").Value = 2688.1875
$ws.Range("J132").Value = 3485.5715
$ws.Range("K132").Value = 8064.5625
$ws.Range("L132").Value = 10456.7145
$ws.Range("M132").Value = -5534.5625
$ws.Range("N132").Value = -15516.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 32778.5
$ws.Range("J56").Value = 32778.5
$ws.Range("L56").Value = 32778.5
$ws.Range("N56").Value = -34160.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 29485
$ws.Range("I58").Value = 28000
$ws.Range("J58").Value = 29980
$ws.Range("K58").Value = 28000
$ws.Range("L58").Value = 29980
$ws.Range("M58").Value = -27692
$ws.Range("N58").Value = -30596

$ws.Range("H61").Value = 17949.5
$ws.Range("J61").Value = 17949.5
$ws.Range("L61").Value = 17949.5
$ws.Range("N61").Value = -18533.5

$ws.Range("H128").Value = 42656.668
$ws.Range("J128").Value = 42656.668
$ws.Range("L128").Value = 42656.668
$ws.Range("N128").Value = -52616.668
